# Apply the "cryptos list" refresh described by the commit:
#   - updated Price / Volume(1h) values for most rows
#   - rows 37/38 swapped (SuiNetwork <-> NEARProtocol) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '56.710.15' },
    @{ Cell = "E2"; Value = '  -0.04%  ' },
    @{ Cell = "D3"; Value = '2.320.41' },
    @{ Cell = "E3"; Value = '  -0.47%  ' },
    @{ Cell = "D5"; Value = '522.98' },
    @{ Cell = "E5"; Value = '  +1.46%  ' },
    @{ Cell = "D6"; Value = '131.93' },
    @{ Cell = "E6"; Value = '  -1.97%  ' },
    @{ Cell = "D7"; Value = '0.995' },
    @{ Cell = "E7"; Value = '  -0.34%  ' },
    @{ Cell = "E8"; Value = '  -0.60%  ' },
    @{ Cell = "D9"; Value = '2.342.53' },
    @{ Cell = "E9"; Value = '  +0.25%  ' },
    @{ Cell = "E10"; Value = '  -1.31%  ' },
    @{ Cell = "E11"; Value = '  +0.45%  ' },
    @{ Cell = "D12"; Value = '5.31' },
    @{ Cell = "E12"; Value = '  -1.26%  ' },
    @{ Cell = "D13"; Value = '0.341' },
    @{ Cell = "E13"; Value = '  -0.28%  ' },
    @{ Cell = "D14"; Value = '23.53' },
    @{ Cell = "E14"; Value = '  -1.56%  ' },
    @{ Cell = "D15"; Value = '2.739.17' },
    @{ Cell = "E15"; Value = '  -0.19%  ' },
    @{ Cell = "D16"; Value = '56.685.11' },
    @{ Cell = "E16"; Value = '  -0.06%  ' },
    @{ Cell = "D17"; Value = '0.0000133' },
    @{ Cell = "E17"; Value = '  -1.29%  ' },
    @{ Cell = "D18"; Value = '2.337.27' },
    @{ Cell = "E18"; Value = '  +0.07%  ' },
    @{ Cell = "D19"; Value = '335.91' },
    @{ Cell = "E19"; Value = '  +2.93%  ' },
    @{ Cell = "D20"; Value = '10.44' },
    @{ Cell = "E20"; Value = '  -0.88%  ' },
    @{ Cell = "D21"; Value = '4.16' },
    @{ Cell = "E21"; Value = '  -1.42%  ' },
    @{ Cell = "D22"; Value = '6.83' },
    @{ Cell = "E22"; Value = '  +3.18%  ' },
    @{ Cell = "D23"; Value = '1.00' },
    @{ Cell = "E23"; Value = '  -0.09%  ' },
    @{ Cell = "D24"; Value = '61.69' },
    @{ Cell = "E24"; Value = '  +1.24%  ' },
    @{ Cell = "D25"; Value = '8.73' },
    @{ Cell = "E26"; Value = '  +0.09%  ' },
    @{ Cell = "D27"; Value = '0.993' },
    @{ Cell = "E27"; Value = '  -0.44%  ' },
    @{ Cell = "E28"; Value = '  +1.80%  ' },
    @{ Cell = "D29"; Value = '169.30' },
    @{ Cell = "E29"; Value = '  -0.45%  ' },
    @{ Cell = "E30"; Value = '  +0.24%  ' },
    @{ Cell = "D31"; Value = '0.0₃0722' },
    @{ Cell = "E31"; Value = '  -2.26%  ' },
    @{ Cell = "D32"; Value = '6.12' },
    @{ Cell = "E32"; Value = '  -1.58%  ' },
    @{ Cell = "D33"; Value = '18.41' },
    @{ Cell = "E33"; Value = '  -0.50%  ' },
    @{ Cell = "D35"; Value = '0.993' },
    @{ Cell = "E35"; Value = '  -0.42%  ' },
    @{ Cell = "E36"; Value = '  -0.77%  ' },
    @{ Cell = "B37"; Value = 'NEARProtocol' },
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = "D37"; Value = '3.98' },
    @{ Cell = "E37"; Value = '  -0.46%  ' },
    @{ Cell = "B38"; Value = 'SuiNetwork' },
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui' },
    @{ Cell = "D38"; Value = '0.900' },
    @{ Cell = "E38"; Value = '  -1.51%  ' },
    @{ Cell = "E39"; Value = '  +1.01%  ' },
    @{ Cell = "D40"; Value = '38.85' },
    @{ Cell = "E40"; Value = '  +1.42%  ' },
    @{ Cell = "D41"; Value = '149.24' },
    @{ Cell = "E41"; Value = '  +4.87%  ' },
    @{ Cell = "E42"; Value = '  -1.09%  ' },
    @{ Cell = "D43"; Value = '287.14' },
    @{ Cell = "E43"; Value = '  +3.44%  ' },
    @{ Cell = "E44"; Value = '  -0.58%  ' },
    @{ Cell = "D45"; Value = '5.10' },
    @{ Cell = "E45"; Value = '  -0.84%  ' },
    @{ Cell = "D46"; Value = '0.0929' },
    @{ Cell = "E46"; Value = '  -0.66%  ' },
    @{ Cell = "D47"; Value = '0.0501' },
    @{ Cell = "E47"; Value = '  -0.81%  ' },
    @{ Cell = "E48"; Value = '  -0.49%  ' },
    @{ Cell = "D49"; Value = '18.58' },
    @{ Cell = "E49"; Value = '  +3.28%  ' },
    @{ Cell = "D50"; Value = '0.0216' },
    @{ Cell = "E50"; Value = '  -1.38%  ' },
    @{ Cell = "E51"; Value = '  -0.64%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text so Excel does not auto-convert numeric-looking
    # strings (e.g. "56.710.15", "0.0000133") into numbers - the
    # source data stores these as plain text.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
